$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row with new columns P (14) and Q (15), copying the
# existing header formatting (bold, bordered, centered) from O1.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("P1").Value2 = 14
$ws.Range("Q1").Value2 = 15

# Rewrite the data block B2:Q25 with the updated simulation results.
# Columns B-K, N-Q change; H, L, M stay zero as before.
$arr = New-Object 'object[,]' 24,16
$arr[0,0] = 2.9966996069532
$arr[0,1] = 1.156917478494364
$arr[0,2] = 0.0114819435364204
$arr[0,3] = 1.43900959886858
$arr[0,4] = 0.9518068714633756
$arr[0,5] = 0.8974281458001059
$arr[0,6] = 0
$arr[0,7] = 0.01989866865855516
$arr[0,8] = 0.51635186465262
$arr[0,9] = 0.3543204419373609
$arr[0,10] = 0
$arr[0,11] = 0
$arr[0,12] = 0
$arr[0,13] = 0
$arr[0,14] = 0.56037868356972
$arr[0,15] = 0
$arr[1,0] = 2.600698824971062
$arr[1,1] = 1.005436972696742
$arr[1,2] = 0.01111827334049664
$arr[1,3] = 1.233109108666511
$arr[1,4] = 0.8458397521930578
$arr[1,5] = 0.7926330485199173
$arr[1,6] = 0
$arr[1,7] = 0.01437554417266096
$arr[1,8] = 0.4740977887499866
$arr[1,9] = 0.3335024851272692
$arr[1,10] = 0
$arr[1,11] = 0
$arr[1,12] = 0
$arr[1,13] = 0
$arr[1,14] = 0.586708392411289
$arr[1,15] = 0
$arr[2,0] = 2.357809705401223
$arr[2,1] = 0.9135254507465334
$arr[2,2] = 0.01089839791584701
$arr[2,3] = 1.109399576780518
$arr[2,4] = 0.7822667999839155
$arr[2,5] = 0.7297689857623908
$arr[2,6] = 0
$arr[2,7] = 0.01138838581976387
$arr[2,8] = 0.44899628231083
$arr[2,9] = 0.3213089467855639
$arr[2,10] = 0
$arr[2,11] = 0
$arr[2,12] = 0
$arr[2,13] = 0
$arr[2,14] = 0.603711179932386
$arr[2,15] = 0
$arr[3,0] = 2.257289175266465
$arr[3,1] = 0.8770664674493958
$arr[3,2] = 0.01083321896932077
$arr[3,3] = 1.059552250100751
$arr[3,4] = 0.755828110680369
$arr[3,5] = 0.7034958258725652
$arr[3,6] = 0
$arr[3,7] = 0.01030537685694233
$arr[3,8] = 0.4384628498483352
$arr[3,9] = 0.3158925565887287
$arr[3,10] = 0
$arr[3,11] = 0
$arr[3,12] = 0
$arr[3,13] = 0
$arr[3,14] = 0.6110488346864891
$arr[3,15] = 0
$arr[4,0] = 2.238690557699499
$arr[4,1] = 0.8719592844028341
$arr[4,2] = 0.010851524078479
$arr[4,3] = 1.051293161254804
$arr[4,4] = 0.7503998153161717
$arr[4,5] = 0.6979505848579919
$arr[4,6] = 0
$arr[4,7] = 0.0101906552431954
$arr[4,8] = 0.4361187851922779
$arr[4,9] = 0.3142909096954725
$arr[4,10] = 0
$arr[4,11] = 0
$arr[4,12] = 0
$arr[4,13] = 0
$arr[4,14] = 0.6125342385192774
$arr[4,15] = 0
$arr[5,0] = 2.351221413989947
$arr[5,1] = 0.9155907124510065
$arr[5,2] = 0.0109772647383064
$arr[5,3] = 1.108685013339397
$arr[5,4] = 0.7790136391765827
$arr[5,5] = 0.7261229249783412
$arr[5,6] = 0
$arr[5,7] = 0.01152041356035483
$arr[5,8] = 0.4471947421380094
$arr[5,9] = 0.3192920316299066
$arr[5,10] = 0
$arr[5,11] = 0
$arr[5,12] = 0
$arr[5,13] = 0
$arr[5,14] = 0.6045201237700049
$arr[5,15] = 0
$arr[6,0] = 2.853047418474148
$arr[6,1] = 1.107851179408414
$arr[6,2] = 0.01146328461280532
$arr[6,3] = 1.367320983468503
$arr[6,4] = 0.9110668007660934
$arr[6,5] = 0.8565633540875837
$arr[6,6] = 0
$arr[6,7] = 0.01806386503106605
$arr[6,8] = 0.4993845233423002
$arr[6,9] = 0.3444242996211315
$arr[6,10] = 0
$arr[6,11] = 0
$arr[6,12] = 0
$arr[6,13] = 0
$arr[6,14] = 0.5702873973686557
$arr[6,15] = 0
$arr[7,0] = 3.851135167235043
$arr[7,1] = 1.491072225367816
$arr[7,2] = 0.01230185355245794
$arr[7,3] = 1.901506410059554
$arr[7,4] = 1.189236995527907
$arr[7,5] = 1.132344697405188
$arr[7,6] = 0
$arr[7,7] = 0.03441531182454671
$arr[7,8] = 0.6124450764323797
$arr[7,9] = 0.4026063668155118
$arr[7,10] = 0
$arr[7,11] = 0
$arr[7,12] = 0
$arr[7,13] = 0
$arr[7,14] = 0.5087716055327576
$arr[7,15] = 0
$arr[8,0] = 4.588253252400023
$arr[8,1] = 1.783851066576517
$arr[8,2] = 0.01298767321211614
$arr[8,3] = 2.318913189592323
$arr[8,4] = 1.403955803021105
$arr[8,5] = 1.345362378587083
$arr[8,6] = 0
$arr[8,7] = 0.04952462130604474
$arr[8,8] = 0.7009804617761404
$arr[8,9] = 0.4487734696399315
$arr[8,10] = 0
$arr[8,11] = 0
$arr[8,12] = 0
$arr[8,13] = 0
$arr[8,14] = 0.4688450600314056
$arr[8,15] = 0
$arr[9,0] = 4.916877338805421
$arr[9,1] = 1.924398032437296
$arr[9,2] = 0.01345003278225931
$arr[9,3] = 2.516351263289621
$arr[9,4] = 1.5001145382045
$arr[9,5] = 1.440186312545023
$arr[9,6] = 0
$arr[9,7] = 0.05743675808424786
$arr[9,8] = 0.7402361707788714
$arr[9,9] = 0.4677557428375536
$arr[9,10] = 0
$arr[9,11] = 0
$arr[9,12] = 0
$arr[9,13] = 0
$arr[9,14] = 0.4533718703952161
$arr[9,15] = 0
$arr[10,0] = 5.046560815953057
$arr[10,1] = 1.9758282874987
$arr[10,2] = 0.01355193980644387
$arr[10,3] = 2.592520034759701
$arr[10,4] = 1.539702766064309
$arr[10,5] = 1.479648989294532
$arr[10,6] = 0
$arr[10,7] = 0.06047415370921616
$arr[10,8] = 0.7568775508433419
$arr[10,9] = 0.4768745593752328
$arr[10,10] = 0
$arr[10,11] = 0
$arr[10,12] = 0
$arr[10,13] = 0
$arr[10,14] = 0.4469193190715615
$arr[10,15] = 0
$arr[11,0] = 5.019582201204855
$arr[11,1] = 1.964262292541946
$arr[11,2] = 0.01351474206355618
$arr[11,3] = 2.576068120168159
$arr[11,4] = 1.531676337381143
$arr[11,5] = 1.471720618612522
$arr[11,6] = 0
$arr[11,7] = 0.05979108018700074
$arr[11,8] = 0.753578777197788
$arr[11,9] = 0.4752489833511575
$arr[11,10] = 0
$arr[11,11] = 0
$arr[11,12] = 0
$arr[11,13] = 0
$arr[11,14] = 0.4481406604404086
$arr[11,15] = 0
$arr[12,0] = 4.92794940393452
$arr[12,1] = 1.928420876137523
$arr[12,2] = 0.01345194226699675
$arr[12,3] = 2.522596076982055
$arr[12,4] = 1.503582878491102
$arr[12,5] = 1.443674425989599
$arr[12,6] = 0
$arr[12,7] = 0.05767384000799591
$arr[12,8] = 0.7417260109064614
$arr[12,9] = 0.4686493086486934
$arr[12,10] = 0
$arr[12,11] = 0
$arr[12,12] = 0
$arr[12,13] = 0
$arr[12,14] = 0.4527695394627997
$arr[12,15] = 0
$arr[13,0] = 4.869973958429114
$arr[13,1] = 1.90745175298855
$arr[13,2] = 0.01344354468561093
$arr[13,3] = 2.489994339329243
$arr[13,4] = 1.485417916526146
$arr[13,5] = 1.425400025716158
$arr[13,6] = 0
$arr[13,7] = 0.05644486655447789
$arr[13,8] = 0.7339180194461221
$arr[13,9] = 0.4639510418060908
$arr[13,10] = 0
$arr[13,11] = 0
$arr[13,12] = 0
$arr[13,13] = 0
$arr[13,14] = 0.455942389288424
$arr[13,15] = 0
$arr[14,0] = 4.550822513688956
$arr[14,1] = 1.78242442396504
$arr[14,2] = 0.01320705948644729
$arr[14,3] = 2.305951275077447
$arr[14,4] = 1.389058365345178
$arr[14,5] = 1.329364046746321
$arr[14,6] = 0
$arr[14,7] = 0.04939113906711068
$arr[14,8] = 0.6935129715593007
$arr[14,9] = 0.4418057613050479
$arr[14,10] = 0
$arr[14,11] = 0
$arr[14,12] = 0
$arr[14,13] = 0
$arr[14,14] = 0.4724708188385875
$arr[14,15] = 0
$arr[15,0] = 4.355799509322082
$arr[15,1] = 1.706511079499933
$arr[15,2] = 0.01306310501454533
$arr[15,3] = 2.195173274338984
$arr[15,4] = 1.330916407205848
$arr[15,5] = 1.271460485380743
$arr[15,6] = 0
$arr[15,7] = 0.04530437946008714
$arr[15,8] = 0.669232528615737
$arr[15,9] = 0.4285847331443406
$arr[15,10] = 0
$arr[15,11] = 0
$arr[15,12] = 0
$arr[15,13] = 0
$arr[15,14] = 0.4829059160543494
$arr[15,15] = 0
$arr[16,0] = 4.248693641943419
$arr[16,1] = 1.660786428931431
$arr[16,2] = 0.01290494357087368
$arr[16,3] = 2.132232323224741
$arr[16,4] = 1.300422333858734
$arr[16,5] = 1.241473679172202
$arr[16,6] = 0
$arr[16,7] = 0.04288369176467466
$arr[16,8] = 0.6569294128516106
$arr[16,9] = 0.4228369601342052
$arr[16,10] = 0
$arr[16,11] = 0
$arr[16,12] = 0
$arr[16,13] = 0
$arr[16,14] = 0.4882335607931054
$arr[16,15] = 0
$arr[17,0] = 4.208754457413136
$arr[17,1] = 1.647141923914148
$arr[17,2] = 0.0129098075043288
$arr[17,3] = 2.11098746267227
$arr[17,4] = 1.288133355521666
$arr[17,5] = 1.229082534927215
$arr[17,6] = 0
$arr[17,7] = 0.04219170857393273
$arr[17,8] = 0.65164378705677
$arr[17,9] = 0.4195751914849595
$arr[17,10] = 0
$arr[17,11] = 0
$arr[17,12] = 0
$arr[17,13] = 0
$arr[17,14] = 0.490655885965424
$arr[17,15] = 0
$arr[18,0] = 4.376912049739985
$arr[18,1] = 1.714381488403205
$arr[18,2] = 0.01307242908016804
$arr[18,3] = 2.20689587110175
$arr[18,4] = 1.337276478371805
$arr[18,5] = 1.277821910933341
$arr[18,6] = 0
$arr[18,7] = 0.0457186249530892
$arr[18,8] = 0.6719159383835347
$arr[18,9] = 0.4301149968099551
$arr[18,10] = 0
$arr[18,11] = 0
$arr[18,12] = 0
$arr[18,13] = 0
$arr[18,14] = 0.4817224852483122
$arr[18,15] = 0
$arr[19,0] = 4.949404756561705
$arr[19,1] = 1.941544639528217
$arr[19,2] = 0.01355565022497629
$arr[19,3] = 2.53817764399021
$arr[19,4] = 1.508869368954322
$arr[19,5] = 1.448542574011668
$arr[19,6] = 0
$arr[19,7] = 0.05843002188812729
$arr[19,8] = 0.7435235463515539
$arr[19,9] = 0.4686333844799435
$arr[19,10] = 0
$arr[19,11] = 0
$arr[19,12] = 0
$arr[19,13] = 0
$arr[19,14] = 0.452303259998061
$arr[19,15] = 0
$arr[20,0] = 5.333383326136811
$arr[20,1] = 2.08959808447554
$arr[20,2] = 0.01376761912348812
$arr[20,3] = 2.762711095844807
$arr[20,4] = 1.628278630367319
$arr[20,5] = 1.568051319409136
$arr[20,6] = 0
$arr[20,7] = 0.06742086343644615
$arr[20,8] = 0.7942755519053719
$arr[20,9] = 0.497590689297617
$arr[20,10] = 0
$arr[20,11] = 0
$arr[20,12] = 0
$arr[20,13] = 0
$arr[20,14] = 0.432950677606911
$arr[20,15] = 0
$arr[21,0] = 5.134726316119441
$arr[21,1] = 2.007136739021348
$arr[21,2] = 0.01355127934790801
$arr[21,3] = 2.642172798309176
$arr[21,4] = 1.5677575477619
$arr[21,5] = 1.507948597520709
$arr[21,6] = 0
$arr[21,7] = 0.06238915503755482
$arr[21,8] = 0.7690299051605791
$arr[21,9] = 0.4843594132296616
$arr[21,10] = 0
$arr[21,11] = 0
$arr[21,12] = 0
$arr[21,13] = 0
$arr[21,14] = 0.4420978258627315
$arr[21,15] = 0
$arr[22,0] = 4.377584059774506
$arr[22,1] = 1.705898659373418
$arr[22,2] = 0.01290830542012955
$arr[22,3] = 2.201734975638885
$arr[22,4] = 1.339949768917904
$arr[22,5] = 1.281253998057849
$arr[22,6] = 0
$arr[22,7] = 0.04525190557496739
$arr[22,8] = 0.6738594778780111
$arr[22,9] = 0.4330976587919295
$arr[22,10] = 0
$arr[22,11] = 0
$arr[22,12] = 0
$arr[22,13] = 0
$arr[22,14] = 0.4806066042566002
$arr[22,15] = 0
$arr[23,0] = 3.570498267296443
$arr[23,1] = 1.390653972596226
$arr[23,2] = 0.01222214960073842
$arr[23,3] = 1.753186020227488
$arr[23,4] = 1.106825924531648
$arr[23,5] = 1.049820793756695
$arr[23,6] = 0
$arr[23,7] = 0.0297451788585974
$arr[23,8] = 0.5778314711576229
$arr[23,9] = 0.3826258538683547
$arr[23,10] = 0
$arr[23,11] = 0
$arr[23,12] = 0
$arr[23,13] = 0
$arr[23,14] = 0.5261099153377984
$arr[23,15] = 0

$ws.Range("B2:Q25").Value2 = $arr
